$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "37.613.28"
$ws.Cells.Item(2, 5).Value = "  +6.18%  "
$ws.Cells.Item(3, 4).Value = "2.052.15"
$ws.Cells.Item(3, 5).Value = "  +3.36%  "
$ws.Cells.Item(4, 5).Value = "  +0.04%  "
$ws.Cells.Item(5, 4).Value = "'253.23"
$ws.Cells.Item(5, 5).Value = "  +5.53%  "
$ws.Cells.Item(6, 5).Value = "  +2.88%  "
$ws.Cells.Item(7, 4).Value = "'65.49"
$ws.Cells.Item(7, 5).Value = "  +18.17%  "
$ws.Cells.Item(8, 5).Value = "  -0.09%  "
$ws.Cells.Item(9, 4).Value = "'0.378"
$ws.Cells.Item(9, 5).Value = "  +7.21%  "
$ws.Cells.Item(10, 4).Value = "'59.37"
$ws.Cells.Item(10, 5).Value = "  +0.40%  "
$ws.Cells.Item(11, 4).Value = "'0.0756"
$ws.Cells.Item(11, 5).Value = "  +4.86%  "
$ws.Cells.Item(12, 5).Value = "  +1.12%  "
$ws.Cells.Item(13, 4).Value = "'0.915"
$ws.Cells.Item(13, 5).Value = "  +3.50%  "
$ws.Cells.Item(14, 4).Value = "'15.21"
$ws.Cells.Item(14, 5).Value = "  +7.30%  "
$ws.Cells.Item(15, 4).Value = "2.345.81"
$ws.Cells.Item(15, 5).Value = "  +3.25%  "
$ws.Cells.Item(16, 4).Value = "'5.63"
$ws.Cells.Item(16, 5).Value = "  +8.86%  "
$ws.Cells.Item(17, 4).Value = "'20.87"
$ws.Cells.Item(17, 5).Value = "  +23.52%  "
$ws.Cells.Item(18, 4).Value = "2.050.82"
$ws.Cells.Item(18, 5).Value = "  +3.36%  "
$ws.Cells.Item(19, 4).Value = "37.458.01"
$ws.Cells.Item(19, 5).Value = "  +5.95%  "
$ws.Cells.Item(20, 4).Value = "'73.28"
$ws.Cells.Item(20, 5).Value = "  +5.60%  "
$ws.Cells.Item(21, 4).Value = "0.0₃0876"
$ws.Cells.Item(21, 5).Value = "  +5.61%  "
$ws.Cells.Item(22, 4).Value = "'5.37"
$ws.Cells.Item(22, 5).Value = "  +8.21%  "
$ws.Cells.Item(23, 4).Value = "'237.76"
$ws.Cells.Item(23, 5).Value = "  +3.11%  "
$ws.Cells.Item(24, 4).Value = "'2.78"
$ws.Cells.Item(24, 5).Value = "  +23.84%  "
$ws.Cells.Item(25, 5).Value = "  +0.24%  "
$ws.Cells.Item(26, 4).Value = "'2.36"
$ws.Cells.Item(26, 5).Value = "  +3.04%  "
$ws.Cells.Item(27, 4).Value = "'9.63"
$ws.Cells.Item(27, 5).Value = "  +6.20%  "
$ws.Cells.Item(28, 4).Value = "'165.95"
$ws.Cells.Item(28, 5).Value = "  +2.34%  "
$ws.Cells.Item(29, 4).Value = "'19.95"
$ws.Cells.Item(29, 5).Value = "  +3.14%  "
$ws.Cells.Item(30, 2).Value = "ImmutableX"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(30, 4).Value = "'1.24"
$ws.Cells.Item(30, 5).Value = "  +10.11%  "
$ws.Cells.Item(31, 2).Value = "Filecoin"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(31, 4).Value = "'5.24"
$ws.Cells.Item(31, 5).Value = "  +10.28%  "
$ws.Cells.Item(32, 2).Value = "Stellar"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(32, 4).Value = "'0.122"
$ws.Cells.Item(32, 5).Value = "  +3.03%  "
$ws.Cells.Item(33, 5).Value = "  +27.30%  "
$ws.Cells.Item(34, 4).Value = "'4.75"
$ws.Cells.Item(34, 5).Value = "  +12.57%  "
$ws.Cells.Item(35, 4).Value = "'0.0618"
$ws.Cells.Item(35, 5).Value = "  +6.20%  "
$ws.Cells.Item(36, 4).Value = "'2.47"
$ws.Cells.Item(36, 5).Value = "  +13.29%  "
$ws.Cells.Item(37, 5).Value = "  -0.01%  "
$ws.Cells.Item(38, 5).Value = "  +0.92%  "
$ws.Cells.Item(39, 4).Value = "'5.99"
$ws.Cells.Item(39, 5).Value = "  +24.47%  "
$ws.Cells.Item(40, 5).Value = "  +19.14%  "
$ws.Cells.Item(41, 4).Value = "'1.24"
$ws.Cells.Item(41, 5).Value = "  +5.54%  "
$ws.Cells.Item(42, 2).Value = "VeChain"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(42, 4).Value = "'0.0220"
$ws.Cells.Item(42, 5).Value = "  +6.57%  "
$ws.Cells.Item(43, 2).Value = "ARBITRUM"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(43, 4).Value = "'1.15"
$ws.Cells.Item(43, 5).Value = "  +7.27%  "
$ws.Cells.Item(44, 4).Value = "'2.74"
$ws.Cells.Item(44, 5).Value = "  +21.75%  "
$ws.Cells.Item(45, 2).Value = "HuobiToken"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(45, 4).Value = "'2.93"
$ws.Cells.Item(45, 5).Value = "  +4.16%  "
$ws.Cells.Item(46, 2).Value = "FraxShare"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(46, 4).Value = "'8.19"
$ws.Cells.Item(46, 5).Value = "  +11.66%  "
$ws.Cells.Item(47, 5).Value = "  +11.94%  "
$ws.Cells.Item(48, 4).Value = "'95.75"
$ws.Cells.Item(48, 5).Value = "  +7.40%  "
$ws.Cells.Item(49, 4).Value = "1.426.07"
$ws.Cells.Item(49, 5).Value = "  +4.48%  "
$ws.Cells.Item(50, 4).Value = "'2.95"
$ws.Cells.Item(50, 5).Value = "  +2.31%  "
$ws.Cells.Item(51, 4).Value = "'47.71"
$ws.Cells.Item(51, 5).Value = "  +6.75%  "
